$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 314 (D314, F314) ---
$ws.Range("D314").Value = 10.1674
$ws.Range("F314").Value = 10.1604

# --- Append new rows 315-317 (monthly FX_IDC:USDMAD OHLC data) ---
# Data for each new row, keyed by target row number.
$newRows = @(
    @{ Row = 315; A = 45170.33333333334; B = "FX_IDC:USDMAD"; C = 10.1598; D = 10.3323; E = 10.112;   F = 10.2833; G = 0 },
    @{ Row = 316; A = 45201.375;         B = "FX_IDC:USDMAD"; C = 10.2833; D = 10.3464; E = 10.1978; F = 10.2975; G = 0 },
    @{ Row = 317; A = 45231.375;         B = "FX_IDC:USDMAD"; C = 10.2958; D = 10.3303; E = 10.1652; F = 10.2043; G = 0 }
)

# Use row 314 as the formatting template (it already has the correct
# date number-format / bold / bordered style for column A, and the
# plain default style for columns B:G), then overwrite the values.
foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A314:G314").Copy()
    $ws.Range("A$row`:G$row").PasteSpecial(-4122)  # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
}

$excel.CutCopyMode = $false
